$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# WebForm User Assignment execution — refresh the phone-number pool
# assigned to the OneY/TwoY rows (column F) and the matched-position
# markers (AM2/AN2). Values are written as text (quote-prefixed) and
# then restored to the Normal style so they keep rendering the same as
# the original plain-text shared-string cells.

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "F2" "9840065345"
Set-TextValue "F3" "9840053220"
Set-TextValue "F4" "9840077394"
Set-TextValue "F5" "9840029072"
Set-TextValue "F6" "9840009216"
Set-TextValue "F7" "9840007879"
Set-TextValue "F8" "9840061451"
Set-TextValue "F9" "9840004707"
Set-TextValue "F10" "9840015495"

Set-TextValue "AM2" "3"
Set-TextValue "AN2" "3"

# Match the refreshed selection left behind by the edit.
$ws.Range("AN2").Select()
